$wb = $excel.ActiveWorkbook

$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

# Update status text (shared by C2/C3 on both sheets via same shared string)
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("C3").Value = "Handed back: in sync with en-US"
$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("C3").Value = "Handed back: in sync with en-US"

# zh-cn: Latest Handback DateTime (H) now populated with the handback timestamp
$zh.Range("H2").Value = "2016-03-18 03:26:22"
$zh.Range("H3").Value = "2016-03-18 03:26:22"

# de-de: Latest Handback DateTime (H) now populated with the handback timestamp
$de.Range("H2").Value = "2016-03-18 03:26:36"
$de.Range("H3").Value = "2016-03-18 03:26:36"

# zh-cn: new Latest Target File (F) / Latest Handback File (G) columns for rows 2 and 3
$zh.Range("F2").Value = "644e43b3-2720-40dc-bbd2-54555b59d7d2.md"
$zh.Range("G2").Value = "644e43b3-2720-40dc-bbd2-54555b59d7d2.8d797c117318588f3ceb7cd358420754a22ce479.zh-cn.xlf"
$zh.Range("F3").Value = "8e59d3a0-f493-47bb-8bdd-835f569b2adb.md"
$zh.Range("G3").Value = "8e59d3a0-f493-47bb-8bdd-835f569b2adb.b5618fcbd3fd8919aa9cf6d27f933d24a94bd620.zh-cn.xlf"

# de-de: new Latest Target File (F) / Latest Handback File (G) columns for rows 2 and 3
$de.Range("F2").Value = "644e43b3-2720-40dc-bbd2-54555b59d7d2.md"
$de.Range("G2").Value = "644e43b3-2720-40dc-bbd2-54555b59d7d2.8d797c117318588f3ceb7cd358420754a22ce479.de-de.xlf"
$de.Range("F3").Value = "8e59d3a0-f493-47bb-8bdd-835f569b2adb.md"
$de.Range("G3").Value = "8e59d3a0-f493-47bb-8bdd-835f569b2adb.b5618fcbd3fd8919aa9cf6d27f933d24a94bd620.de-de.xlf"

# Hyperlinks for the new F/G cells, duplicating the targets used by A/D in the same row
$zh.Range("F2").Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/07068a46b2f0daf2a5c50a48e95e2c1a6afa9c1a/e2e/644e43b3-2720-40dc-bbd2-54555b59d7d2.md", "", "", "644e43b3-2720-40dc-bbd2-54555b59d7d2.md") | Out-Null
$zh.Range("G2").Hyperlinks.Add($zh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b4d4f3990d145978a818bc05552a5cc62b600776/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/644e43b3-2720-40dc-bbd2-54555b59d7d2.8d797c117318588f3ceb7cd358420754a22ce479.zh-cn.xlf", "", "", "644e43b3-2720-40dc-bbd2-54555b59d7d2.8d797c117318588f3ceb7cd358420754a22ce479.zh-cn.xlf") | Out-Null
$zh.Range("F3").Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/07068a46b2f0daf2a5c50a48e95e2c1a6afa9c1a/e2e/8e59d3a0-f493-47bb-8bdd-835f569b2adb.md", "", "", "8e59d3a0-f493-47bb-8bdd-835f569b2adb.md") | Out-Null
$zh.Range("G3").Hyperlinks.Add($zh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b4d4f3990d145978a818bc05552a5cc62b600776/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/8e59d3a0-f493-47bb-8bdd-835f569b2adb.b5618fcbd3fd8919aa9cf6d27f933d24a94bd620.zh-cn.xlf", "", "", "8e59d3a0-f493-47bb-8bdd-835f569b2adb.b5618fcbd3fd8919aa9cf6d27f933d24a94bd620.zh-cn.xlf") | Out-Null

$de.Range("F2").Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/07068a46b2f0daf2a5c50a48e95e2c1a6afa9c1a/e2e/644e43b3-2720-40dc-bbd2-54555b59d7d2.md", "", "", "644e43b3-2720-40dc-bbd2-54555b59d7d2.md") | Out-Null
$de.Range("G2").Hyperlinks.Add($de.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fea6ace382183a1505fb14023d93927474d70861/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/644e43b3-2720-40dc-bbd2-54555b59d7d2.8d797c117318588f3ceb7cd358420754a22ce479.de-de.xlf", "", "", "644e43b3-2720-40dc-bbd2-54555b59d7d2.8d797c117318588f3ceb7cd358420754a22ce479.de-de.xlf") | Out-Null
$de.Range("F3").Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/07068a46b2f0daf2a5c50a48e95e2c1a6afa9c1a/e2e/8e59d3a0-f493-47bb-8bdd-835f569b2adb.md", "", "", "8e59d3a0-f493-47bb-8bdd-835f569b2adb.md") | Out-Null
$de.Range("G3").Hyperlinks.Add($de.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fea6ace382183a1505fb14023d93927474d70861/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/8e59d3a0-f493-47bb-8bdd-835f569b2adb.b5618fcbd3fd8919aa9cf6d27f933d24a94bd620.de-de.xlf", "", "", "8e59d3a0-f493-47bb-8bdd-835f569b2adb.b5618fcbd3fd8919aa9cf6d27f933d24a94bd620.de-de.xlf") | Out-Null

# Match the same "HyperLink" font appearance (underline, blue) used by the existing
# A/B/D hyperlink cells on these rows, instead of the engine's auto-generated default
$newLinkCells = @($zh.Range("F2"), $zh.Range("G2"), $zh.Range("F3"), $zh.Range("G3"), $de.Range("F2"), $de.Range("G2"), $de.Range("F3"), $de.Range("G3"))
foreach ($cell in $newLinkCells) {
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 11
    $cell.Font.Underline = 2
    $cell.Font.Color = 15570276
}
